$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the existing header cell H1 into I1 and J1
$ws.Range("H1").Copy() | Out-Null
$ws.Range("I1:J1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# Set header values
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data values for columns I (I0) and J (IF), rows 2-16
$values = @(
    @(7, 8),
    @(8, 8),
    @(9, 9),
    @(5, 5),
    @(7, 7),
    @(7, 7),
    @(7, 7),
    @(7, 7),
    @(7, 7),
    @(7, 7),
    @(7, 7),
    @(6, 6),
    @(5, 5),
    @(5, 5),
    @(5, 5)
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $values[$i][0]
    $ws.Cells.Item($row, 10).Value = $values[$i][1]
}
